# Rename PascalCase column/field names to camelCase on the "Toinen versio"
# sheet, to match the naming convention already used on the second sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Toinen versio")

$ws.Range("C3").Value  = "eventId"
$ws.Range("G3").Value  = "ticketId"
$ws.Range("K3").Value  = "ticketTypeId"

$ws.Range("C4").Value  = "eventName"
$ws.Range("G4").Value  = "eventId"
$ws.Range("K4").Value  = "ticketType"

$ws.Range("C5").Value  = "date"
$ws.Range("G5").Value  = "ticketTypeId"
$ws.Range("K5").Value  = "price"

$ws.Range("C6").Value  = "ticketCount"

$ws.Range("C7").Value  = "venueId"

$ws.Range("C8").Value  = "description"

$ws.Range("C11").Value = "venueId"

$ws.Range("C12").Value = "place"

$ws.Range("C13").Value = "streetAddress"

$ws.Range("C14").Value = "postalCode"

$ws.Range("C15").Value = "cityId"

$ws.Range("C18").Value = "cityId"

$ws.Range("C19").Value = "city"

$ws.Range("G20").Value = "username"
$ws.Range("K20").Value = "roleId"

$ws.Range("G21").Value = "password"
$ws.Range("K21").Value = "role"

$ws.Range("C22").Value = "postalCode"
$ws.Range("G22").Value = "roleId"

$ws.Range("C23").Value = "postOffice"
